$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row: Date (D), Volumen (J), Precio minimo (K), Precio maximo (L),
# Precio promedio ponderado (M), Precio $/Kg (P)
$data = @{
    2  = @{ D = 44280; J = 55; K = 4000; L = 4000; M = 4000; P = 4000 }
    3  = @{ D = 44390; J = 55; K = 6000; L = 6000; M = 6000; P = 6000 }
    4  = @{ D = 44312; J = 50; K = 4000; L = 4000; M = 4000; P = 4000 }
    5  = @{ D = 44259; J = 30; K = 4000; L = 4000; M = 4000; P = 4000 }
    6  = @{ D = 44365; J = 55; K = 5000; L = 5000; M = 5000; P = 5000 }
    7  = @{ D = 44291; J = 35; K = 4000; L = 4000; M = 4000; P = 4000 }
    8  = @{ D = 44316; J = 20; K = 4000; L = 4000; M = 4000; P = 4000 }
    9  = @{ D = 44301; J = 40; K = 3000; L = 3000; M = 3000; P = 3000 }
    10 = @{ D = 44176; J = 10; K = 4000; L = 4000; M = 4000; P = 4000 }
    11 = @{ D = 44315; J = 40; K = 4000; L = 4000; M = 4000; P = 4000 }
    12 = @{ D = 44313; J = 20; K = 4000; L = 4000; M = 4000; P = 4000 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("P$row").Value = $vals.P
}
